$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly-off timestamp previously written for row 11
# (same instant, just re-serialized with the task's native precision).
$ws.Range("A11").Value2 = 45862.87522489583

# Append the new reading captured by the scheduled task (row 12).
$ws.Range("A12").Value2 = 45862.9168948495
$ws.Range("A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B12").Value2 = 2025
$ws.Range("C12").Value2 = 30
$ws.Range("D12").Value2 = 19.33
$ws.Range("E12").Value2 = 74.44
$ws.Range("F12").Value2 = 87.89
$ws.Range("G12").Value2 = 13.54
$ws.Range("H12").Value = "ESE"
$ws.Range("I12").Value2 = 0
$ws.Range("J12").Value = "22:00:19"
